# Weekly fruit/vegetable data update: a new daily observation is inserted
# at the top of the data block (row 203), pushing the existing rows 203-263
# down by one (to 204-264). This grows the sheet's used range from
# A1:R263 to A1:R264.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 203:263 down to 204:264, leaving row 203 blank for the new record.
$ws.Rows("203:203").Insert()

# Populate the newly inserted row 203 with the new observation.
$ws.Range("A203").Value = 6
$ws.Range("B203").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C203").Value = "Metropolitana"
$ws.Range("D203").Value = 44841
$ws.Range("E203").Value = 13
$ws.Range("F203").Value = 100112022
$ws.Range("G203").Value = "Arveja Verde"
$ws.Range("H203").Value = "Perfection"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 400
$ws.Range("K203").Value = 27000
$ws.Range("L203").Value = 28000
$ws.Range("M203").Value = 27425
$ws.Range("N203").Value = "$/malla 25 kilos"
$ws.Range("O203").Value = "Provincia de Huasco"
$ws.Range("P203").Value = 1097
$ws.Range("Q203").Value = 25
$ws.Range("R203").Value = "Hortaliza"
